$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.144.33"
$ws.Range("E2").Value = "  -4.35%  "
$ws.Range("D3").Value = "1.652.51"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "215.99"
$ws.Range("E5").Value = "  -3.77%  "
$ws.Range("D6").Value = "0.5109"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("E9").Value = "  -3.69%  "
$ws.Range("D10").Value = "19.97"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "1.655.55"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "4.279"
$ws.Range("E13").Value = "  -4.87%  "
$ws.Range("D14").Value = "1.880.85"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").Value = "0.5518"
$ws.Range("E15").Value = "  -5.34%  "
$ws.Range("D16").Value = "0.0₅8028"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("E17").Value = "  -5.78%  "
$ws.Range("D18").Value = "26.153.40"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "210.66"
$ws.Range("E20").Value = "  -4.74%  "
$ws.Range("E21").Value = "  -4.68%  "
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").Value = "6.036"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D25").Value = "143.57"
$ws.Range("D26").Value = "1.736"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("D27").Value = "0.1182"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").Value = "6.980"
$ws.Range("E28").Value = "  -3.51%  "
$ws.Range("D29").Value = "15.82"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("D30").Value = "0.05109"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").Value = "3.341"
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("E33").Value = "  -6.08%  "
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").Value = "2.362"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").Value = "1.168.70"
$ws.Range("E38").Value = "  +5.04%  "
$ws.Range("D39").Value = "0.5689"
$ws.Range("D40").Value = "0.01587"
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "0.8302"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").Value = "5.653"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").Value = "100.38"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "1.790.80"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "0.4551"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "55.64"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "7.877"
$ws.Range("E51").Value = "  -2.96%  "
